$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-ambiguous cell updates (coin names, links, percentages, and
#     price strings that are not parseable as a single plain number) ---
$ws.Range("D2").Value = "43.032.50"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "2.289.02"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  -3.02%  "
$ws.Range("E6").Value = "  +4.11%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("E10").Value = "  -4.19%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("E14").Value = "  +21.29%  "
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "2.631.92"
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").Value = "2.289.26"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").Value = "42.983.71"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("E21").Value = "  -4.27%  "
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("E23").Value = "  +3.12%  "
$ws.Range("E24").Value = "  +4.06%  "
$ws.Range("E25").Value = "  +10.44%  "
$ws.Range("E26").Value = "  -4.28%  "
$ws.Range("E27").Value = "  -4.09%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E29").Value = "  -2.87%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("E31").Value = "  +4.81%  "
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("E33").Value = "  -2.84%  "
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("E36").Value = "  +6.59%  "
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("E38").Value = "  -9.46%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("E41").Value = "  -3.45%  "
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E44").Value = "  +5.34%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E46").Value = "  -8.19%  "
$ws.Range("E47").Value = "  +1.92%  "
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("E49").Value = "  +3.43%  "
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("E51").Value = "  -0.78%  "

# --- Numeric-looking price strings: must stay as literal text (e.g. "1.00", "0.0983")
#     Force text format so Excel does not coerce them to numbers, then reset the
#     cell style back to Normal so no stray number-format style is left behind. ---
$numericTextCells = @{
    "D4" = "1.00"
    "D5" = "113.00"
    "D6" = "310.00"
    "D7" = "0.634"
    "D9" = "0.617"
    "D10" = "44.56"
    "D11" = "0.0929"
    "D12" = "55.09"
    "D13" = "8.86"
    "D14" = "1.08"
    "D16" = "15.53"
    "D21" = "7.21"
    "D22" = "76.56"
    "D24" = "2.47"
    "D25" = "257.34"
    "D26" = "9.02"
    "D27" = "11.77"
    "D28" = "1.00"
    "D29" = "39.05"
    "D30" = "2.24"
    "D31" = "22.35"
    "D32" = "173.89"
    "D34" = "0.0902"
    "D35" = "5.75"
    "D36" = "5.09"
    "D37" = "0.130"
    "D38" = "4.18"
    "D39" = "0.0377"
    "D42" = "72.51"
    "D43" = "0.232"
    "D44" = "1.42"
    "D45" = "1.00"
    "D47" = "5.72"
    "D48" = "108.46"
    "D49" = "8.93"
    "D50" = "1.31"
    "D51" = "0.0983"
}
foreach ($addr in $numericTextCells.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $numericTextCells.Keys) {
    $ws.Range($addr).Value = $numericTextCells[$addr]
}
foreach ($addr in $numericTextCells.Keys) {
    $ws.Range($addr).Style = "Normal"
}
